# Applies the "Add descriptions titles" edit:
#  - Metadata sheet: Title value, Description value, and Date value
#  - Elements sheet: Short / Definition / Mapping: RIM Mapping for the
#    root "Extension" row

$wb = $excel.ActiveWorkbook

$meta = $wb.Worksheets.Item("Metadata")
$meta.Range("B5").Value = "DMI Organization Location"
$meta.Range("B8").Value = "2026-02-25T08:15:31+00:00"
$meta.Range("B12").Value = "Extension créée dans ce volet pour représenter le lieu de l'organisation."

$elements = $wb.Worksheets.Item("Elements")
$elements.Range("L2").Value = "DMI Organization Location"
$elements.Range("M2").Value = "Extension créée dans ce volet pour représenter le lieu de l'organisation."
$elements.Range("AK2").Value = ""
